# Corrected up-to-date Burn Down chart data.
# The underlying "Burn Down Chart Table" worksheet drives the "Burn Down
# Chart" chartsheet via formulas (row 30 = remaining work per day, which
# the chart plots). Updating the raw input cells below recalculates the
# AC (per-task remaining) column and the row 30/31 roll-ups automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn Down Chart Table")

$ws.Range("D2").Value = 1
$ws.Range("D5").Value = 0.5
$ws.Range("D6").Value = 0.5
$ws.Range("F8").Value = 3
$ws.Range("E9").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("G11").Value = 2
$ws.Range("D12").Value = 1
$ws.Range("H12").Value = 0.5
$ws.Range("M14").Value = 0.5
$ws.Range("K19").Value = 0.5
$ws.Range("M19").Value = 0.5
$ws.Range("L20").Value = 0.5

# Restore the default top-left cell and move the selection, matching the
# author's latest view of the table.
$ws.Range("L16").Select()
